# "commit and add commands are also added"
# Updates the Git practice sheet: fixes a couple of labels, bumps the
# practice-times counter, and appends two new sections ("Add" and
# "Commit") below the existing "Branches" section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- small existing-data tweaks -------------------------------------------

# Practice times for cloning went from 3 to 10
$ws.Range("D2").Value = 10

# Wording fixes on the "Branches" section
$ws.Range("A8").Value = "To list branches that are currently in working  directory"
$ws.Range("A9").Value = "To list all the branches (local & remote)"

# C11 / C12 pick up the bold "command" style used by the rest of the column
$ws.Range("C11").Font.Bold = $true
$ws.Range("C12").Font.Bold = $true

# --- new "Commit" / "Add" section headers -----------------------------------

$ws.Range("A18").Value = "Commit"
$ws.Range("A14").Value = "Add"

# --- new "Add" section rows (15-16) ----------------------------------------

$ws.Range("A15").Value = "To add the work done in the working directory"
$ws.Range("B15").Value = "use add"
$ws.Range("C15").Value = "git add 'File name'"
$ws.Range("C15").Font.Bold = $true

$ws.Range("A16").Value = "To add all the files"
$ws.Range("C16").Value = "git add ."
$ws.Range("C16").Font.Bold = $true
$ws.Range("B16").Value = "add "

# --- new "Commit" section rows (19-20) --------------------------------------

$ws.Range("A19").Value = "To commit the work done"
$ws.Range("B19").Value = "commit -m"
$ws.Range("C19").Value = "git commit -m 'Message to be given'"
$ws.Range("C19").Font.Bold = $true

$ws.Range("C20").Value = "git commit -m 'Subject' -m 'Description'"
$ws.Range("C20").Font.Bold = $true

# --- sheet view / print setup ----------------------------------------------

[void]$ws.Range("C19:C20").Select()

$ws.PageSetup.Orientation = 1
